# Update "想去人数" (interested-people count) figures in sheets "展览" (F column)
# and "全部类型" (F column, rows offset by +1 because it aggregates all rows).
# This mirrors a scraped-data refresh (gh-pages data update).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new value
$wsExhibit.Range("F4").Value  = 1443
$wsExhibit.Range("F5").Value  = 333
$wsExhibit.Range("F7").Value  = 10819
$wsExhibit.Range("F10").Value = 300
$wsExhibit.Range("F13").Value = 12137
$wsExhibit.Range("F14").Value = 12597

# Sheet "全部类型": same events, rows shifted down by 1
$wsAll.Range("F5").Value  = 1443
$wsAll.Range("F6").Value  = 333
$wsAll.Range("F8").Value  = 10819
$wsAll.Range("F11").Value = 300
$wsAll.Range("F14").Value = 12137
$wsAll.Range("F15").Value = 12597
